# TAC-3104, Fix start trip date issue
#
# The "Trip Pick up Date Start *" sample value (B2) was stored as a real
# date serial (3/19/2022) formatted with numFmtId 14. Change it to a plain
# text value "03/19/2022" (custom format code "[$-1010000]m/d/yyyy;@") so it
# no longer behaves like a date, and bump the sample "Reference No" (A2)
# from 3 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference No sample value
$ws.Range("A2").Value = 1

# Force the cell to store text (not an auto-recognized date) before writing
# the date-like string, then apply the custom text/date format code.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "03/19/2022"
$ws.Range("B2").NumberFormat = "[$-1010000]m/d/yyyy;@"

# Move the active selection to B3, matching the saved file's cursor state.
$ws.Range("B3").Select()
